# This script recomputes the transient-response values in the circuit-analysis
# result sheets after the user equation (forcing function v(t) = t^3 applied at
# node 1) was added to the simulation code. It writes the freshly computed
# branch/node voltages and branch currents (previously ~0, i.e. numerical noise
# from an unforced/trivial circuit) into the result tables for t = 0.1 ... 1.0 s.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Branches_Voltages")
$ws.Range("B3").Value = 0.001
$ws.Range("C3").Value = -0.0009824278318694955
$ws.Range("D3").Value = 0.00001757216813050468
$ws.Range("E3").Value = -0.000009932095030285254
$ws.Range("F3").Value = -0.000007640073100219423
$ws.Range("B4").Value = 0.008000000000000002
$ws.Range("C4").Value = -0.007844163559570896
$ws.Range("D4").Value = 0.0001558364404291058
$ws.Range("E4").Value = -0.00007944486195533348
$ws.Range("F4").Value = -0.0000763915784737723
$ws.Range("B5").Value = 0.02700000000000001
$ws.Range("C5").Value = -0.02638891150794136
$ws.Range("D5").Value = 0.0006110884920586503
$ws.Range("E5").Value = -0.0002676787503065252
$ws.Range("F5").Value = -0.0003434097417521251
$ws.Range("B6").Value = 0.06400000000000002
$ws.Range("C6").Value = -0.06233185230461617
$ws.Range("D6").Value = 0.001668147695383843
$ws.Range("E6").Value = -0.0006323828649629163
$ws.Range("F6").Value = -0.001035764830420926
$ws.Range("B7").Value = 0.125
$ws.Range("C7").Value = -0.1212996897185551
$ws.Range("D7").Value = 0.003700310281444933
$ws.Range("E7").Value = -0.001231098357116649
$ws.Range("F7").Value = -0.002469211924328284
$ws.Range("B8").Value = 0.2160000000000001
$ws.Range("C8").Value = -0.2088336968083816
$ws.Range("D8").Value = 0.007166303191618545
$ws.Range("E8").Value = -0.002119617517548127
$ws.Range("F8").Value = -0.005046685674070419
$ws.Range("B9").Value = 0.3430000000000001
$ws.Range("C9").Value = -0.3303889817942656
$ws.Range("D9").Value = 0.01261101820573446
$ws.Range("E9").Value = -0.003353919466789186
$ws.Range("F9").Value = -0.009257098738945273
$ws.Range("B10").Value = 0.5120000000000001
$ws.Range("C10").Value = -0.4913381910815169
$ws.Range("D10").Value = 0.02066180891848324
$ws.Range("E10").Value = -0.004987914681134856
$ws.Range("F10").Value = -0.01567389423734838
$ws.Range("B11").Value = 0.7290000000000001
$ws.Range("C11").Value = -0.6969700852614578
$ws.Range("D11").Value = 0.03202991473854228
$ws.Range("E11").Value = -0.00707604868279009
$ws.Range("F11").Value = -0.02495386605575219
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = -0.9524937969742547
$ws.Range("D12").Value = 0.04750620302574527
$ws.Range("E12").Value = -0.009670430164435181
$ws.Range("F12").Value = -0.03783577286131009

$ws = $wb.Worksheets.Item("Nodes_Voltages")
$ws.Range("B3").Value = -0.001
$ws.Range("C3").Value = -0.00001757216813050468
$ws.Range("D3").Value = -0.000007640073100219423
$ws.Range("B4").Value = -0.008000000000000002
$ws.Range("C4").Value = -0.0001558364404291058
$ws.Range("D4").Value = -0.0000763915784737723
$ws.Range("B5").Value = -0.02700000000000001
$ws.Range("C5").Value = -0.0006110884920586503
$ws.Range("D5").Value = -0.0003434097417521251
$ws.Range("B6").Value = -0.06400000000000002
$ws.Range("C6").Value = -0.001668147695383843
$ws.Range("D6").Value = -0.001035764830420926
$ws.Range("B7").Value = -0.125
$ws.Range("C7").Value = -0.003700310281444933
$ws.Range("D7").Value = -0.002469211924328284
$ws.Range("B8").Value = -0.2160000000000001
$ws.Range("C8").Value = -0.007166303191618545
$ws.Range("D8").Value = -0.005046685674070419
$ws.Range("B9").Value = -0.3430000000000001
$ws.Range("C9").Value = -0.01261101820573446
$ws.Range("D9").Value = -0.009257098738945273
$ws.Range("B10").Value = -0.5120000000000001
$ws.Range("C10").Value = -0.02066180891848324
$ws.Range("D10").Value = -0.01567389423734838
$ws.Range("B11").Value = -0.7290000000000001
$ws.Range("C11").Value = -0.03202991473854228
$ws.Range("D11").Value = -0.02495386605575219
$ws.Range("B12").Value = -1
$ws.Range("C12").Value = -0.04750620302574527
$ws.Range("D12").Value = -0.03783577286131009

$ws = $wb.Worksheets.Item("Branches_Currents")
$ws.Range("B3").Value = -0.00007675217436480435
$ws.Range("C3").Value = -0.00007675217436480433
$ws.Range("D3").Value = 0.0000003514433626100935
$ws.Range("E3").Value = -0.00007640073100219427
$ws.Range("F3").Value = -0.00007640073100219423
$ws.Range("B4").Value = -0.000766329626821085
$ws.Range("C4").Value = -0.0007663296268210849
$ws.Range("D4").Value = 0.000002413842083361928
$ws.Range("E4").Value = -0.0007639157847377231
$ws.Range("F4").Value = -0.0007639157847377229
$ws.Range("B5").Value = -0.00344078861647048
$ws.Range("C5").Value = -0.00344078861647048
$ws.Range("D5").Value = 0.000006691198949228962
$ws.Range("E5").Value = -0.003434097417521252
$ws.Range("F5").Value = -0.003434097417521252
$ws.Range("B6").Value = -0.01037209828932654
$ws.Range("C6").Value = -0.01037209828932654
$ws.Range("D6").Value = 0.00001444998511727488
$ws.Range("E6").Value = -0.01035764830420926
$ws.Range("F6").Value = -0.01035764830420926
$ws.Range("B7").Value = -0.02471831250988679
$ws.Range("C7").Value = -0.02471831250988679
$ws.Range("D7").Value = 0.00002619326660394694
$ws.Range("E7").Value = -0.02469211924328284
$ws.Range("F7").Value = -0.02469211924328284
$ws.Range("B8").Value = -0.05050998333230372
$ws.Range("C8").Value = -0.05050998333230372
$ws.Range("D8").Value = 0.0000431265915995253
$ws.Range("E8").Value = -0.0504668567407042
$ws.Range("F8").Value = -0.05046685674070419
$ws.Range("B9").Value = -0.09263675509813554
$ws.Range("C9").Value = -0.09263675509813553
$ws.Range("D9").Value = 0.00006576770868279297
$ws.Range("E9").Value = -0.09257098738945276
$ws.Range("F9").Value = -0.09257098738945273
$ws.Range("B10").Value = -0.156834190479056
$ws.Range("C10").Value = -0.156834190479056
$ws.Range("D10").Value = 0.00009524810557218265
$ws.Range("E10").Value = -0.1567389423734838
$ws.Range("F10").Value = -0.1567389423734838
$ws.Range("B11").Value = -0.2496707745683509
$ws.Range("C11").Value = -0.2496707745683509
$ws.Range("D11").Value = 0.0001321140108289982
$ws.Range("E11").Value = -0.2495386605575219
$ws.Range("F11").Value = -0.2495386605575219
$ws.Range("B12").Value = -0.3785351403680159
$ws.Range("C12").Value = -0.3785351403680159
$ws.Range("D12").Value = 0.0001774117549150616
$ws.Range("E12").Value = -0.3783577286131009
$ws.Range("F12").Value = -0.3783577286131009
